# Delete row 5 ("  " / "  ") on the NegativeLoginTest sheet -- this shifts the
# old row 6 (S12345 / "  ") up to become the new row 5, shrinking the used
# range from A1:B6 to A1:B5 and dropping the shared-string ref count by 2.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("NegativeLoginTest")

$ws1.Rows("5:5").Delete()

# Make NegativeLoginTest the active/selected sheet (was CampusBankData before)
# with its selection anchored on the new last cell, B5.
$ws1.Activate()
$ws1.Range("B5").Select()
